# Weekly update: insert two new data rows (780:781) for the latest week,
# shifting the existing historical rows down by two (the two rows that
# overflow the previous used range are preserved at the new end, rows
# 816:817).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 780-781; everything below shifts down by 2
# (old row 814/815 content lands on new rows 816/817 automatically).
$ws.Rows("780:781").Insert()

# Populate the newly inserted row 780 (Primera calidad, new week).
$ws.Cells.Item(780, 1).Value2  = 6
$ws.Cells.Item(780, 2).Value2  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(780, 3).Value2  = "Metropolitana"
$ws.Cells.Item(780, 4).Value2  = 44753
$ws.Cells.Item(780, 5).Value2  = 13
$ws.Cells.Item(780, 6).Value2  = 100112023
$ws.Cells.Item(780, 7).Value2  = "Brócoli"
$ws.Cells.Item(780, 8).Value2  = "Sin especificar"
$ws.Cells.Item(780, 9).Value2  = "Primera"
$ws.Cells.Item(780, 10).Value2 = 6700
$ws.Cells.Item(780, 11).Value2 = 900
$ws.Cells.Item(780, 12).Value2 = 1000
$ws.Cells.Item(780, 13).Value2 = 949
$ws.Cells.Item(780, 14).Value2 = "$/unidad"
$ws.Cells.Item(780, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(780, 16).Value2 = 949
$ws.Cells.Item(780, 17).Value2 = 1
$ws.Cells.Item(780, 18).Value2 = "Hortaliza"

# Populate the newly inserted row 781 (Segunda calidad, new week).
$ws.Cells.Item(781, 1).Value2  = 6
$ws.Cells.Item(781, 2).Value2  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(781, 3).Value2  = "Metropolitana"
$ws.Cells.Item(781, 4).Value2  = 44753
$ws.Cells.Item(781, 5).Value2  = 13
$ws.Cells.Item(781, 6).Value2  = 100112023
$ws.Cells.Item(781, 7).Value2  = "Brócoli"
$ws.Cells.Item(781, 8).Value2  = "Sin especificar"
$ws.Cells.Item(781, 9).Value2  = "Segunda"
$ws.Cells.Item(781, 10).Value2 = 2600
$ws.Cells.Item(781, 11).Value2 = 700
$ws.Cells.Item(781, 12).Value2 = 700
$ws.Cells.Item(781, 13).Value2 = 700
$ws.Cells.Item(781, 14).Value2 = "$/unidad"
$ws.Cells.Item(781, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(781, 16).Value2 = 700
$ws.Cells.Item(781, 17).Value2 = 1
$ws.Cells.Item(781, 18).Value2 = "Hortaliza"
